# The presentation's Design/Theme was changed from the "Integral" (Red
# Violet) theme to the default "Office Theme" (Office colour scheme).
#
# In the underlying OOXML this shows up as the 12 theme colours
# (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) stored in
# ppt/theme/theme1.xml's <a:clrScheme> changing from the "Red Violet"
# palette to the standard "Office" palette (the font scheme and format
# scheme are already identical between the Office and Integral themes,
# so no change is needed there).
#
# PowerPoint exposes the active theme's colours through
# Slide.ThemeColorScheme (a 12-entry collection of RGBColor objects,
# in the fixed order dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
# Because every slide shares the single slide master/theme, updating
# the scheme through any one slide updates it for the whole deck.

function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$colorScheme = $s.ThemeColorScheme

# Target palette: the standard Office theme colours.
$officeColors = @(
    (RGB 0x00 0x00 0x00),  #  1 dk1
    (RGB 0xFF 0xFF 0xFF),  #  2 lt1
    (RGB 0x44 0x54 0x6A),  #  3 dk2
    (RGB 0xE7 0xE6 0xE6),  #  4 lt2
    (RGB 0x5B 0x9B 0xD5),  #  5 accent1
    (RGB 0xED 0x7D 0x31),  #  6 accent2
    (RGB 0xA5 0xA5 0xA5),  #  7 accent3
    (RGB 0xFF 0xC0 0x00),  #  8 accent4
    (RGB 0x44 0x72 0xC4),  #  9 accent5
    (RGB 0x70 0xAD 0x47),  # 10 accent6
    (RGB 0x05 0x63 0xC1),  # 11 hlink
    (RGB 0x95 0x4F 0x72)   # 12 folHlink
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
